{"js": "// Find the run of text that starts the target paragraph segment.\nconst body = context.document.body;\nconst results = body.search(\"En las b\u00fasquedas manuales no fue mucho lo que se encontr\u00f3\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\nconst match = results.items[0];\n\n// Extend the matched run's text in place (preserves the run's existing\n// formatting: highlight=white, rtl=0) and append the new sentences\n// separated by line breaks (vertical-tab \\u000b => <w:br/> line break),\n// ending with one extra break that takes the place of the old italic\n// \"i=1\" run's lone break (that run gets folded away below).\nconst newTail =\n  \", incluso en google dorks no se encontraron archivos con informaci\u00f3n sensible.\" +\n  \"\\u000b\\u000b\" +\n  \"creemos que tienen un nivel de seguridad aceptable por lo menos en cuanto al manejo de los archivos\u2026\" +\n  \"\\u000b\\u000b\" +\n  \"respecto a las conexiones...\" +\n  \"\\u000b\";\nmatch.insertText(newTail, Word.InsertLocation.end);\nawait context.sync();\n\n// Locate the tail of the paragraph (the three remaining line breaks):\n//   1) the break we just typed above (keep)\n//   2) the old now-orphaned italic run's break (remove)\n//   3) the old blue (#0000FF) run's break (keep, untouched)\nconst results2 = body.search(\"respecto a las conexiones...\", { matchCase: true });\nresults2.load(\"text\");\nawait context.sync();\nconst tailAnchor = results2.items[0];\n\nconst paragraph = tailAnchor.paragraphs.getFirst();\nconst paragraphEnd = paragraph.getRange(Word.RangeLocation.end);\nconst afterAnchor = tailAnchor.getRange(Word.RangeLocation.after);\nconst tailRange = afterAnchor.expandTo(paragraphEnd);\ntailRange.load(\"text\");\nawait context.sync();\n\nconst breakRanges = tailRange.getTextRanges([\"\\u000b\"], false);\nbreakRanges.load(\"text,items\");\nawait context.sync();\nfor (const part of breakRanges.items) {\n  part.font.load(\"italic\");\n}\nawait context.sync();\n\n// Delete the single leftover break that belonged to the removed italic run.\nconst orphan = breakRanges.items.find((part) => part.font.italic === true);\nif (orphan) {\n  orphan.delete();\n  await context.sync();\n}\n", "ps1": "# Locate the run of text that starts the target paragraph segment.\n$d = $word.ActiveDocument\n$rng = $d.Content\n$f = $rng.Find\n$f.ClearFormatting()\n$f.Text = \"En las b\u00fasquedas manuales no fue mucho lo que se encontr\u00f3\"\n$f.MatchCase = $true\n$f.Forward = $true\n$ok = $f.Execute()\nif (-not $ok) {\n    throw \"Target text not found\"\n}\n\n# Remember the enclosing paragraph before we mutate anything (Find has\n# narrowed $rng down to just the matched text).\n$targetParagraph = $rng.Paragraphs(1)\n\n# Extend the matched run's text in place (this preserves the run's existing\n# formatting: highlight=white, rtl=0) by appending the new sentences,\n# separated by line breaks (Chr(11) => a <w:br/> line break). The text ends\n# with one extra break that takes over for the old italic \"i=1\" run's lone\n# break; that now-redundant run is folded away below.\n$rng.Collapse(0)  # wdCollapseEnd\n$br = [char]11\n$newText = \", incluso en google dorks no se encontraron archivos con informaci\u00f3n sensible.\" + $br + $br + `\n    \"creemos que tienen un nivel de seguridad aceptable por lo menos en cuanto al manejo de los archivos\u2026\" + $br + $br + `\n    \"respecto a las conexiones...\" + $br\n$rng.InsertAfter($newText)\n\n# The paragraph now ends with three consecutive line breaks:\n#   1) the break we just typed above (keep)\n#   2) the old now-orphaned italic run's break (remove)\n#   3) the old blue (#0000FF) run's break (keep, untouched)\n$paraRange = $targetParagraph.Range\n$tailStart = $paraRange.End - 4\n$tailRange = $d.Range($tailStart, $paraRange.End)\n\nfor ($k = 1; $k -le $tailRange.Characters.Count; $k++) {\n    $ch = $tailRange.Characters($k)\n    if (([int][char]$ch.Text[0]) -eq 11 -and $ch.Font.Italic) {\n        $ch.Delete()\n        break\n    }\n}\n"}
